$wb = $excel.ActiveWorkbook
$admin = $wb.Worksheets.Item("Admin")

# --- Rebuild the "Admin" sheet content to its final layout ---
$admin.Cells.Clear()

$admin.Range("A1").Value = "Admin Modules"

$admin.Range("A2").Value = "Dashboard:"
$admin.Range("E2").Value = "Type"

$admin.Range("B3").Value = "API Stats"
$admin.Range("E3").Value = "Table"

$admin.Range("C4").Value = "API Succesful Requests"
$admin.Range("C5").Value = "API Errors"
$admin.Range("C6").Value = "Server Exceptions (HTTP 500)"
$admin.Range("C7").Value = "Bets"
$admin.Range("C8").Value = "Suggestions"
$admin.Range("C9").Value = "Admins"
$admin.Range("C10").Value = "Users"
$admin.Range("C11").Value = "Books"

$admin.Range("B13").Value = "Book Stats"
$admin.Range("C14").Value = "Book1"
$admin.Range("D15").Value = "API Succesful Requests"
$admin.Range("D16").Value = "API Errors"
$admin.Range("D17").Value = "Bets"
$admin.Range("D18").Value = "Suggestions"

$admin.Range("C19").Value = "Book2"
$admin.Range("D20").Value = "API Succesful Requests"
$admin.Range("D21").Value = "API Errors"
$admin.Range("D22").Value = "Bets"
$admin.Range("D23").Value = "Suggestions"

$admin.Range("B25").Value = "API Activity Chart last 24 hrs"
$admin.Range("E25").Value = "line area chart"
$admin.Range("B26").Value = "API Activity Chart last 7 days"
$admin.Range("E26").Value = "line area chart"
$admin.Range("B27").Value = "API Activity Chart this month"
$admin.Range("E27").Value = "line area chart"

$admin.Range("A29").Value = "Admins:"
$admin.Range("B30").Value = "Page to add or remove admins"

$admin.Range("A32").Value = "Users:"
$admin.Range("B33").Value = "Page to add or remove users"
$admin.Range("B34").Value = "Page to add a book to a certain user"

$admin.Range("A36").Value = "Books:"
$admin.Range("B37").Value = "Page to add, edit/config or disable books"

$admin.Range("A39").Value = "Reports:"
$admin.Range("B40").Value = "API Stats"
$admin.Range("B41").Value = "Book Stats"
$admin.Range("B42").Value = "Financials"
$admin.Range("B43").Value = "Suggestions Tracker"

$admin.Range("D15").Select()

# --- New "Admin TODO" sheet ---
$todo = $wb.Worksheets.Add($null, $admin)
$todo.Name = "Admin TODO"

$todo.Range("A1").Value = "Todo:"
$todo.Range("A2").Value = "Dashboard"
$todo.Range("A3").Value = "Admins"
$todo.Range("A4").Value = "Users"
$todo.Range("A5").Value = "Books"
$todo.Range("A6").Value = "Reports"

$todo.Range("B1").Select()
$todo.Activate()
